$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two trailing rows (75 and 76) so the sheet shrinks from A1:C76 to A1:C74,
# matching the updated dataset (2020 Q1/Q2 replace the old 2019 Q3/Q4 rows, and several
# interior rows change from quarter-to-quarter GDP-style shocks to the corrected MP-surprise shocks).
$ws.Rows.Item(76).Delete() | Out-Null
$ws.Rows.Item(75).Delete() | Out-Null

# row, year, quarter, agg_shock (corrected MP-surprise based shock values)
$data = @(
    @(2, 2001, 2, -0.004739260866827267),
    @(3, 2001, 3, -0.01194664645579203),
    @(4, 2002, 4, -0.002258249291871672),
    @(5, 2003, 1, 0.0203566002794262),
    @(6, 2003, 2, 0.01691399345752377),
    @(7, 2003, 3, 0.0406971664029268),
    @(8, 2003, 4, 0.06675379425142157),
    @(9, 2004, 1, -0.03486890844848055),
    @(10, 2004, 2, -0.0191714759926876),
    @(11, 2004, 3, -0.007569254409319139),
    @(12, 2004, 4, -0.0008051471035429339),
    @(13, 2005, 1, -0.0036266134203475),
    @(14, 2005, 2, -0.0346107602739647),
    @(15, 2005, 3, -0.0003647771170994296),
    @(16, 2005, 4, 0.001823403185341064),
    @(17, 2006, 1, -0.000399189012530506),
    @(18, 2006, 2, 0.03064551655092188),
    @(19, 2006, 3, 0.001709962340397238),
    @(20, 2006, 4, -0.001488206372120969),
    @(21, 2007, 1, 0.001305829990442957),
    @(22, 2007, 2, -0.03033331254723405),
    @(23, 2007, 3, 0.006751987552642134),
    @(24, 2007, 4, -0.0103195447518296),
    @(25, 2008, 1, -0.01803337069679335),
    @(26, 2008, 2, -0.0240482526936212),
    @(27, 2008, 3, -0.03383209414086327),
    @(28, 2008, 4, -0.207437509019274),
    @(29, 2009, 1, 0.1685179598429348),
    @(30, 2009, 2, 0.0387720554634715),
    @(31, 2009, 3, -0.03507804968142587),
    @(32, 2009, 4, -0.009462151537709303),
    @(33, 2010, 1, -0.06360301516382566),
    @(34, 2010, 2, 0.02945869019485453),
    @(35, 2010, 3, 0.003194846775386333),
    @(36, 2010, 4, -0.02171917110487806),
    @(37, 2011, 1, -0.0075062351385305),
    @(38, 2011, 2, -0.02756997217903178),
    @(39, 2011, 3, -0.03047701866701014),
    @(40, 2011, 4, -0.03066979709712514),
    @(41, 2012, 1, 0.06669769684404064),
    @(42, 2012, 2, -0.009506673238740363),
    @(43, 2012, 3, 0.023364962987077),
    @(44, 2012, 4, -0.0007041072120317996),
    @(45, 2013, 1, 0.01872107263379565),
    @(46, 2013, 2, -0.0177893906595021),
    @(47, 2013, 3, 0.02279953986076947),
    @(48, 2013, 4, 0.0087593643538472),
    @(49, 2014, 1, 0.0075950358713036),
    @(50, 2014, 2, 0.00132975434098792),
    @(51, 2014, 3, 0.003775470945794947),
    @(52, 2014, 4, 0.004858046556991167),
    @(53, 2015, 1, 0.003004337637710901),
    @(54, 2015, 2, 0.03802878612740767),
    @(55, 2015, 3, 0.01738058398764043),
    @(56, 2015, 4, 0.03135040227546287),
    @(57, 2016, 1, 0.0094228259850608),
    @(58, 2016, 2, -0.0267298967152486),
    @(59, 2016, 3, -0.01647688698675413),
    @(60, 2016, 4, 0.004636066349316954),
    @(61, 2017, 1, 0.002724627479826017),
    @(62, 2017, 2, -0.0001183915858551757),
    @(63, 2017, 3, -0.001074620177498969),
    @(64, 2017, 4, 0.0000494962201844543),
    @(65, 2018, 1, -0.0003872654054307455),
    @(66, 2018, 2, -0.0001184668024195385),
    @(67, 2018, 3, 0.0004934021675000028),
    @(68, 2018, 4, -0.0001949615029347847),
    @(69, 2019, 1, 0.004722577775968351),
    @(70, 2019, 2, 0.04284479502328777),
    @(71, 2019, 3, 0.003173945619475534),
    @(72, 2019, 4, -0.009729173632727767),
    @(73, 2020, 1, 0.03669913426772903),
    @(74, 2020, 2, 0.0013710353099937)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $data[$i][0]
    $year = $data[$i][1]
    $quarter = $data[$i][2]
    $shock = $data[$i][3]
    $ws.Cells.Item($r, 1).Value = $year
    $ws.Cells.Item($r, 2).Value = $quarter
    $ws.Cells.Item($r, 3).Value = $shock
}

Write-Output "Updated $($data.Length) rows; dimension now A1:C74"
